# Insert a new data row at row 618 (pushing the existing rows 618-714 down
# to 619-715) and populate it with a new price-record for "Papa" (potato)
# at Terminal Hortofrutícola Agro Chillán, matching the rest of the sheet's
# constant columns (A, B, C, E, F, G, Q, R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("618:618").Insert()

$ws.Range("A618").Value = 7
$ws.Range("B618").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C618").Value = "Ñuble"
$ws.Range("D618").Value = 45154
$ws.Range("E618").Value = 16
$ws.Range("F618").Value = 100114001
$ws.Range("G618").Value = "Papa"
$ws.Range("H618").Value = "Asterix"
$ws.Range("I618").Value = "2a (guarda)"
$ws.Range("J618").Value = 200
$ws.Range("K618").Value = 19000
$ws.Range("L618").Value = 19000
$ws.Range("M618").Value = 19000
$ws.Range("N618").Value = "`$/saco 25 kilos"
$ws.Range("O618").Value = "Región de Los Lagos"
$ws.Range("P618").Value = 760
$ws.Range("Q618").Value = 25
$ws.Range("R618").Value = "Hortaliza"
